$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 372.25
$ws.Range("I12").Value = 305.8
$ws.Range("J12").Value = 483
$ws.Range("K12").Value = 305.8
$ws.Range("L12").Value = 483
$ws.Range("M12").Value = -135.8
$ws.Range("N12").Value = -823
$ws.Range("H43").Value = 3099.5
$ws.Range("I43").Value = 2632.6667
$ws.Range("J43").Value = 4500
$ws.Range("K43").Value = 2632.6667
$ws.Range("L43").Value = 4500
$ws.Range("M43").Value = -2563.6667
$ws.Range("N43").Value = -4638
$ws.Range("H70").Value = 10000000
$ws.Range("J70").Value = 10000000
$ws.Range("L70").Value = 30000000
$ws.Range("N70").Value = -30000540
$ws.Range("H73").Value = 10000000
$ws.Range("J73").Value = 10000000
$ws.Range("L73").Value = 30000000
$ws.Range("N73").Value = -30001872
$ws.Range("H86").Value = 433.33334
$ws.Range("I86").Value = 433.33334
$ws.Range("K86").Value = 433.33334
$ws.Range("M86").Value = 689.66666
$ws.Range("H89").Value = 433.33334
$ws.Range("I89").Value = 433.33334
$ws.Range("K89").Value = 2166.6667
$ws.Range("M89").Value = 3449.3333
$ws.Range("H124").Value = 184000
$ws.Range("J124").Value = 184000
$ws.Range("L124").Value = 184000
$ws.Range("N124").Value = -193820
$ws.Range("H135").Value = 1812.1428
$ws.Range("I135").Value = 1812.1428
$ws.Range("K135").Value = 16309.2852
$ws.Range("M135").Value = -13774.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3314.6296
$ws.Range("I32").Value = 2819.84
$ws.Range("K32").Value = 2819.84
$ws.Range("M32").Value = -2532.84
$ws.Range("H45").Value = 1292.5
$ws.Range("I45").Value = 1292.5
$ws.Range("K45").Value = 1292.5
$ws.Range("M45").Value = -915.5
$ws.Range("H61").Value = 1823.8823
$ws.Range("I61").Value = 1528.0714
$ws.Range("J61").Value = 3204.3333
$ws.Range("K61").Value = 1528.0714
$ws.Range("L61").Value = 3204.3333
$ws.Range("M61").Value = -1316.0714
$ws.Range("N61").Value = -3628.3333
$ws.Range("H110").Value = 2073.6667
$ws.Range("I110").Value = 1444
$ws.Range("K110").Value = 1444
$ws.Range("M110").Value = 601
$ws.Range("H132").Value = 2049.0952
$ws.Range("I132").Value = 1039.7037
$ws.Range("K132").Value = 3119.1111
$ws.Range("M132").Value = -589.1111000000001
$ws.Range("H136").Value = 1823.8823
$ws.Range("I136").Value = 1528.0714
$ws.Range("J136").Value = 3204.3333
$ws.Range("K136").Value = 4584.2142
$ws.Range("L136").Value = 9612.999899999999
$ws.Range("M136").Value = -2034.2142
$ws.Range("N136").Value = -14712.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H107").Value = 2051.6667
$ws.Range("I107").Value = 1859.6666
$ws.Range("K107").Value = 1859.6666
$ws.Range("M107").Value = 60.33339999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1934
$ws.Range("I58").Value = 1192.6666
$ws.Range("K58").Value = 1192.6666
$ws.Range("M58").Value = -989.6666
$ws.Range("H105").Value = 1496.5
$ws.Range("I105").Value = 1496.5
$ws.Range("K105").Value = 1496.5
$ws.Range("M105").Value = 250.5
$ws.Range("H107").Value = 1139.3077
$ws.Range("I107").Value = 502
$ws.Range("K107").Value = 502
$ws.Range("M107").Value = 1418
$ws.Range("H122").Value = 2418
$ws.Range("I122").Value = 1881.4445
$ws.Range("J122").Value = 3625.25
$ws.Range("K122").Value = 5644.333500000001
$ws.Range("L122").Value = 10875.75
$ws.Range("M122").Value = -3194.333500000001
$ws.Range("N122").Value = -15775.75
$ws.Range("H132").Value = 2238
$ws.Range("I132").Value = 1357.5
$ws.Range("K132").Value = 4072.5
$ws.Range("M132").Value = -1542.5
$ws.Range("H134").Value = 2081.625
$ws.Range("I134").Value = 2207.5715
$ws.Range("K134").Value = 6622.7145
$ws.Range("M134").Value = -4087.7145
$ws.Range("H136").Value = 1934
$ws.Range("I136").Value = 1192.6666
$ws.Range("K136").Value = 3577.9998
$ws.Range("M136").Value = -1027.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 497.33334
$ws.Range("I8").Value = 497.33334
$ws.Range("K8").Value = 1492.00002
$ws.Range("M8").Value = -1353.00002
$ws.Range("H129").Value = 1549.3334
$ws.Range("I129").Value = 1549.3334
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 4648.0002
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 351.9997999999996
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 3666.5
$ws.Range("I132").Value = 3194
$ws.Range("J132").Value = 4139
$ws.Range("K132").Value = 28746
$ws.Range("L132").Value = 37251
$ws.Range("M132").Value = -26216
$ws.Range("N132").Value = -42311
$ws.Range("H139").Value = 1986.6666
$ws.Range("I139").Value = 1986.6666
$ws.Range("K139").Value = 5959.9998
$ws.Range("M139").Value = -819.9997999999996
$ws.Range("H140").Value = 3249.5
$ws.Range("I140").Value = 3249.5
$ws.Range("K140").Value = 9748.5
$ws.Range("M140").Value = -4568.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3317.5557
$ws.Range("I132").Value = 2030.5
$ws.Range("J132").Value = 3685.2856
$ws.Range("K132").Value = 6091.5
$ws.Range("L132").Value = 11055.8568
$ws.Range("M132").Value = -3561.5
$ws.Range("N132").Value = -16115.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3232.5
$ws.Range("I7").Value = 3194.2856
$ws.Range("K7").Value = 3194.2856
$ws.Range("M7").Value = -3082.2856
$ws.Range("H22").Value = 520.3333
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 195
$ws.Range("H27").Value = 520.3333
$ws.Range("I27").Value = 100
$ws.Range("K27").Value = 100
$ws.Range("M27").Value = 7
$ws.Range("H40").Value = 2312.75
$ws.Range("I40").Value = 2312.75
$ws.Range("K40").Value = 2312.75
$ws.Range("M40").Value = -2176.75
$ws.Range("H122").Value = 1652
$ws.Range("I122").Value = 1652
$ws.Range("K122").Value = 4956
$ws.Range("M122").Value = -2506
$ws.Range("H126").Value = 3232.5
$ws.Range("I126").Value = 3194.2856
$ws.Range("K126").Value = 9582.856800000001
$ws.Range("M126").Value = -7112.856800000001
$ws.Range("H132").Value = 4424.5
$ws.Range("I132").Value = 3875.75
$ws.Range("J132").Value = 4698.875
$ws.Range("K132").Value = 11627.25
$ws.Range("L132").Value = 14096.625
$ws.Range("M132").Value = -9097.25
$ws.Range("N132").Value = -19156.625
$ws.Range("H136").Value = 11427264
$ws.Range("I136").Value = 15996869
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 47990607
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -47988057
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1289.8572
$ws.Range("I113").Value = 1226.4
$ws.Range("J113").Value = 1448.5
$ws.Range("K113").Value = 3679.2
$ws.Range("L113").Value = 4345.5
$ws.Range("M113").Value = -1509.2
$ws.Range("N113").Value = -8685.5
$ws.Range("H132").Value = 3116.16
$ws.Range("I132").Value = 2640.8572
$ws.Range("J132").Value = 3721.0908
$ws.Range("K132").Value = 7922.571599999999
$ws.Range("L132").Value = 11163.2724
$ws.Range("M132").Value = -5392.571599999999
$ws.Range("N132").Value = -16223.2724
$ws.Range("H136").Value = 2172.5789
$ws.Range("I136").Value = 1572.75
$ws.Range("J136").Value = 3200.8572
$ws.Range("K136").Value = 4718.25
$ws.Range("L136").Value = 9602.571599999999
$ws.Range("M136").Value = -2168.25
$ws.Range("N136").Value = -14702.5716
